$d = $word.ActiveDocument
$d.Content.Find.Execute("{{ contact.name.first }}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{{ contact.name.full() }}", 2)
